# Add June 2021 (1..30) rows (rows 450-479) to every sheet of the workbook.
# Row 450 carries real data (C, D[, E]) for each sheet; rows 451-479 only
# carry the date in column A (data not yet available for those days).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Per-sheet data for the first new day (row 450, date serial 44348 = 2021-06-01)
# ---------------------------------------------------------------------------
$sheetData = @(
    @{ name = "Nuovi casi";        cValue = 21; styleC = $false },
    @{ name = "Deceduti";          cValue = 0;  styleC = $false },
    @{ name = "Dimessi   Guariti"; cValue = 76; styleC = $false },
    @{ name = "Ricoveri";          cValue = 38; styleC = $true; hasE = $true },
    @{ name = "Terapia";           cValue = 4;  styleC = $true }
)

foreach ($info in $sheetData) {
    $ws = $wb.Worksheets.Item($info.name)

    # Reference cells (the last existing data row, 449) used to clone styles.
    $refA = $ws.Range("A449")
    $refC = $ws.Range("C449")
    $refD = $ws.Range("D449")

    # Fill in the date column (A) for all 30 new days (rows 450-479).
    for ($r = 450; $r -le 479; $r++) {
        $ws.Cells.Item($r, 1).Value = 44347 + ($r - 449)
    }
    $ws.Range("A450:A479").NumberFormat = $refA.NumberFormat

    # Row 450 - the only new row that already has data.
    $ws.Range("C450").Value = $info.cValue
    if ($info.styleC) {
        # C column on this sheet uses a font-only style (no custom number format).
        $ws.Range("C450").Font.Color = $refC.Font.Color
    }

    $ws.Range("D450").Formula = "=AVERAGE(C444:C450)"
    $ws.Range("D450").Font.Color = $refD.Font.Color
    $ws.Range("D450").NumberFormat = $refD.NumberFormat

    if ($info.hasE) {
        $ws.Range("E450").Formula = "=C450-C449"
    }
}

# ---------------------------------------------------------------------------
# View state: the active selection on every sheet moves from row 449 to the
# new last row (450), and "Nuovi casi" becomes the active/selected tab
# (instead of "Terapia").
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("Deceduti").Range("A450:D450").Select()
$wb.Worksheets.Item("Dimessi   Guariti").Range("A450:D450").Select()
$wb.Worksheets.Item("Ricoveri").Range("A450:D450").Select()
$wb.Worksheets.Item("Terapia").Range("A450:A479").Select()
$wb.Worksheets.Item("Nuovi casi").Range("A450:D450").Select()
$wb.Worksheets.Item("Nuovi casi").Activate()
